$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035308180912854
$ws.Range("D2").Value = 1.053007218568018
$ws.Range("E2").Value = 1.034119133837381
$ws.Range("F2").Value = 1.057544371387219
$ws.Range("I2").Value = 1.040809748337796
$ws.Range("J2").Value = 1.040422665769014
$ws.Range("K2").Value = 1.055754377434827
$ws.Range("L2").Value = 1.036919524428249
$ws.Range("M2").Value = 1.0602790502247
$ws.Range("N2").Value = 1.041900185292554
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038300985665554
$ws.Range("D3").Value = 1.054258277299736
$ws.Range("E3").Value = 1.036735000896053
$ws.Range("F3").Value = 1.059398956129997
$ws.Range("I3").Value = 1.04122405548643
$ws.Range("J3").Value = 1.043048671635841
$ws.Range("K3").Value = 1.056817272340243
$ws.Range("L3").Value = 1.039339705122441
$ws.Range("M3").Value = 1.061944851116097
$ws.Range("N3").Value = 1.044529920388919
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0402242690657
$ws.Range("D4").Value = 1.0550605577591
$ws.Range("E4").Value = 1.038415663561509
$ws.Range("F4").Value = 1.060589195123923
$ws.Range("I4").Value = 1.041487130357193
$ws.Range("J4").Value = 1.044734983504674
$ws.Range("K4").Value = 1.057497257032361
$ws.Range("L4").Value = 1.0408935346045
$ws.Range("M4").Value = 1.063012526407238
$ws.Range("N4").Value = 1.046218627014029
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041029731697383
$ws.Range("D5").Value = 1.055396128757229
$ws.Range("E5").Value = 1.039119423855581
$ws.Range("F5").Value = 1.061087267927612
$ws.Range("I5").Value = 1.041596538775869
$ws.Range("J5").Value = 1.045440902016094
$ws.Range("K5").Value = 1.057781284353437
$ws.Range("L5").Value = 1.041543920078927
$ws.Range("M5").Value = 1.063458973775612
$ws.Range("N5").Value = 1.046925548010704
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.0411647940275
$ws.Range("D6").Value = 1.055452373183365
$ws.Range("E6").Value = 1.039237426901681
$ws.Range("F6").Value = 1.061170762681487
$ws.Range("I6").Value = 1.041614839566631
$ws.Range("J6").Value = 1.045559254769775
$ws.Range("K6").Value = 1.057828866686567
$ws.Range("L6").Value = 1.041652957978201
$ws.Range("M6").Value = 1.063533794490368
$ws.Range("N6").Value = 1.047044068838871
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040235043703803
$ws.Range("D7").Value = 1.055065048351842
$ws.Range("E7").Value = 1.038425078102716
$ws.Range("F7").Value = 1.060595859392127
$ws.Range("I7").Value = 1.041488596932151
$ws.Range("J7").Value = 1.04474442773368
$ws.Range("K7").Value = 1.057501059410212
$ws.Range("L7").Value = 1.040902236166323
$ws.Range("M7").Value = 1.063018501248088
$ws.Range("N7").Value = 1.046228084654923
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036322413392361
$ws.Range("D8").Value = 1.053431533772584
$ws.Range("E8").Value = 1.035005705599821
$ws.Range("F8").Value = 1.058173192500805
$ws.Range("I8").Value = 1.040950809383764
$ws.Range("J8").Value = 1.041312853571544
$ws.Range("K8").Value = 1.05611521325845
$ws.Range("L8").Value = 1.037740004553571
$ws.Range("M8").Value = 1.060844154204966
$ws.Range("N8").Value = 1.042791637263869
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029322210173106
$ws.Range("D9").Value = 1.050496478861981
$ws.Range("E9").Value = 1.028885053480067
$ws.Range("F9").Value = 1.05382711082098
$ws.Range("I9").Value = 1.039964281910023
$ws.Range("J9").Value = 1.035163652899761
$ws.Range("K9").Value = 1.053612481286353
$ws.Range("L9").Value = 1.032071017164081
$ws.Range("M9").Value = 1.056932610287759
$ws.Range("N9").Value = 1.036633704021809
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024578552038499
$ws.Range("D10").Value = 1.048500063491633
$ws.Range("E10").Value = 1.024735505805264
$ws.Range("F10").Value = 1.050875115179025
$ws.Range("I10").Value = 1.039279724109462
$ws.Range("J10").Value = 1.030990234416847
$ws.Range("K10").Value = 1.051901587359487
$ws.Range("L10").Value = 1.028221839845995
$ws.Range("M10").Value = 1.054268395227992
$ws.Range("N10").Value = 1.032454358805952
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022504935507702
$ws.Range("D11").Value = 1.047625791698519
$ws.Range("E11").Value = 1.022921159158263
$ws.Range("F11").Value = 1.049583270137849
$ws.Range("I11").Value = 1.038976759488789
$ws.Range("J11").Value = 1.029164380264286
$ws.Range("K11").Value = 1.051150324896351
$ws.Range("L11").Value = 1.02653743511357
$ws.Range("M11").Value = 1.053100738036456
$ws.Range("N11").Value = 1.030625911731066
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0217316432666
$ws.Range("D12").Value = 1.047299540408052
$ws.Range("E12").Value = 1.022244489756791
$ws.Range("F12").Value = 1.049101319433492
$ws.Range("I12").Value = 1.038863226315164
$ws.Range("J12").Value = 1.028483258364874
$ws.Range("K12").Value = 1.050869672347557
$ws.Range("L12").Value = 1.02590901851592
$ws.Range("M12").Value = 1.052664855190924
$ws.Range("N12").Value = 1.02994382256043
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021897657407417
$ws.Range("D13").Value = 1.047369591205353
$ws.Range("E13").Value = 1.022389763304102
$ws.Range("F13").Value = 1.049204795547303
$ws.Range("I13").Value = 1.038887625003883
$ws.Range("J13").Value = 1.02862949508258
$ws.Range("K13").Value = 1.050929946220886
$ws.Range("L13").Value = 1.026043942203086
$ws.Range("M13").Value = 1.052758452313614
$ws.Range("N13").Value = 1.03009026695105
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022441078016063
$ws.Range("D14").Value = 1.047598854624046
$ws.Range("E14").Value = 1.022865281952342
$ws.Range("F14").Value = 1.04954347515723
$ws.Range("I14").Value = 1.03896739526194
$ws.Range("J14").Value = 1.029108138659291
$ws.Range("K14").Value = 1.051127158905696
$ws.Range("L14").Value = 1.026485546745423
$ws.Range("M14").Value = 1.053064752269792
$ws.Range("N14").Value = 1.030569590256541
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0227754885211
$ws.Range("D15").Value = 1.04773991055088
$ws.Range("E15").Value = 1.023157898491507
$ws.Range("F15").Value = 1.049751866581745
$ws.Range("I15").Value = 1.039016411563825
$ws.Range("J15").Value = 1.029402656905132
$ws.Range("K15").Value = 1.051248455111487
$ws.Range("L15").Value = 1.026757266050174
$ws.Range("M15").Value = 1.053253185509169
$ws.Range("N15").Value = 1.030864526752073
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024715748686446
$ws.Range("D16").Value = 1.048557876463641
$ws.Range("E16").Value = 1.02485553935953
$ws.Range("F16").Value = 1.050960559079007
$ws.Range("I16").Value = 1.039299691669297
$ws.Range("J16").Value = 1.031111006829368
$ws.Range("K16").Value = 1.051951223580295
$ws.Range("L16").Value = 1.028333247442891
$ws.Range("M16").Value = 1.05434558849013
$ws.Range("N16").Value = 1.032575302729151
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025927498739002
$ws.Range("D17").Value = 1.049068314679317
$ws.Range("E17").Value = 1.025915651344912
$ws.Range("F17").Value = 1.051715056102476
$ws.Range("I17").Value = 1.039475622604294
$ws.Range("J17").Value = 1.032177521264348
$ws.Range("K17").Value = 1.052389235469893
$ws.Range("L17").Value = 1.029317016677194
$ws.Range("M17").Value = 1.055027028672791
$ws.Range("N17").Value = 1.033643331736959
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026632410800473
$ws.Range("D18").Value = 1.049365100400377
$ws.Range("E18").Value = 1.026532308939843
$ws.Range("F18").Value = 1.052153832079622
$ws.Range("I18").Value = 1.039577609536506
$ws.Range("J18").Value = 1.032797800435936
$ws.Range("K18").Value = 1.052643715911215
$ws.Range("L18").Value = 1.029889132894364
$ws.Range("M18").Value = 1.055423150840445
$ws.Range("N18").Value = 1.034264491776119
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026872451996388
$ws.Range("D19").Value = 1.049466137509716
$ws.Range("E19").Value = 1.026742289991961
$ws.Range("F19").Value = 1.05230322295715
$ws.Range("I19").Value = 1.03961227793242
$ws.Range("J19").Value = 1.033008997271418
$ws.Range("K19").Value = 1.052730317760136
$ws.Range("L19").Value = 1.030083924321965
$ws.Range("M19").Value = 1.055557990843041
$ws.Range("N19").Value = 1.034475988535332
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025797684846299
$ws.Range("D20").Value = 1.049013647399824
$ws.Range("E20").Value = 1.02580208666107
$ws.Range("F20").Value = 1.051634241553451
$ws.Range("I20").Value = 1.039456812213989
$ws.Range("J20").Value = 1.032063281347375
$ws.Range("K20").Value = 1.052342345070567
$ws.Range("L20").Value = 1.02921164404986
$ws.Range("M20").Value = 1.054954056653146
$ws.Range("N20").Value = 1.033528929586201
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022281139582557
$ws.Range("D21").Value = 1.047531384161833
$ws.Range("E21").Value = 1.022725330035361
$ws.Range("F21").Value = 1.04944380094317
$ws.Range("I21").Value = 1.038943932599551
$ws.Range("J21").Value = 1.028967271462649
$ws.Range("K21").Value = 1.051069129126295
$ws.Range("L21").Value = 1.026355582019238
$ws.Range("M21").Value = 1.052974614721155
$ws.Range("N21").Value = 1.030428523012321
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020052390720751
$ws.Range("D22").Value = 1.046590685104626
$ws.Range("E22").Value = 1.020774944323822
$ws.Range("F22").Value = 1.048054393357297
$ws.Range("I22").Value = 1.038615679271974
$ws.Range("J22").Value = 1.027003748923101
$ws.Range("K22").Value = 1.050259332439305
$ws.Range("L22").Value = 1.024543880593966
$ws.Range("M22").Value = 1.051717518520937
$ws.Range("N22").Value = 1.028462212045548
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021235615774498
$ws.Range("D23").Value = 1.04709020763751
$ws.Range("E23").Value = 1.021810423067727
$ws.Range("F23").Value = 1.048792119615634
$ws.Range("I23").Value = 1.038790246139256
$ws.Range("J23").Value = 1.02804629025998
$ws.Range("K23").Value = 1.050689511337437
$ws.Range("L23").Value = 1.02550584570021
$ws.Range("M23").Value = 1.052385136506321
$ws.Range("N23").Value = 1.029506233910707
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025856347962978
$ws.Range("D24").Value = 1.04903835210645
$ws.Range("E24").Value = 1.025853406862857
$ws.Range("F24").Value = 1.051670762212575
$ws.Range("I24").Value = 1.039465313767528
$ws.Range("J24").Value = 1.032114907006575
$ws.Range("K24").Value = 1.052363535925912
$ws.Range("L24").Value = 1.02925926265109
$ws.Range("M24").Value = 1.054987033737579
$ws.Range("N24").Value = 1.03358062855976
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031145040078
$ws.Range("D25").Value = 1.051262135423201
$ws.Range("E25").Value = 1.030479191592467
$ws.Range("F25").Value = 1.054960090630867
$ws.Range("I25").Value = 1.04022400412005
$ws.Range("J25").Value = 1.036766016689544
$ws.Range("K25").Value = 1.054266851522654
$ws.Range("L25").Value = 1.033548537093665
$ws.Range("M25").Value = 1.057953598108535
$ws.Range("N25").Value = 1.038238343351968

